# Adds summary statistic rows/formulas to the freelancer results sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15-18: headline summary statistics with bold, larger, vertically
# centered labels/values.
$ws.Range("A15").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(N2:N31)"

$ws.Range("A16").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B16").Formula = "=AVERAGE(Z2:Z31)"

$ws.Range("A17").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B17").Formula = "=MIN(N2:N31)"

$ws.Range("A18").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B18").Formula = "=MAX(Z2:Z31)"

$summaryRange = $ws.Range("B15:B18")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108

$ws.Rows("15:15").RowHeight = 15.6
$ws.Rows("16:16").RowHeight = 15.6
$ws.Rows("17:17").RowHeight = 15.6
$ws.Rows("18:18").RowHeight = 15.6

# Row 12: average of the J column (k fraction) over the 10 instances.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Page setup: letter-ish A4/9-size portrait print, matching the re-saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection left on J12 after the edits, as in the saved workbook.
$ws.Range("J12").Select() | Out-Null
